$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 575; this shifts the
# existing rows 575-593 down to 576-594 (dimension becomes A1:R594).
$ws.Rows(575).Insert()

$ws.Range("A575").Value = 9
$ws.Range("B575").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C575").Value = "Metropolitana"
$ws.Range("D575").Value = 45075
$ws.Range("E575").Value = 13
$ws.Range("F575").Value = 100112052
$ws.Range("G575").Value = "Albahaca"
$ws.Range("H575").Value = "Sin especificar"
$ws.Range("I575").Value = "Primera"
$ws.Range("J575").Value = 520
$ws.Range("K575").Value = 2500
$ws.Range("L575").Value = 3000
$ws.Range("M575").Value = 2750
$ws.Range("N575").Value = "$/paquete"
$ws.Range("O575").Value = "Región de Arica y Parinacota"
$ws.Range("P575").Value = 2750
$ws.Range("Q575").Value = 1
$ws.Range("R575").Value = "Hortaliza"
